$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend formatting for the newly added rows (21-27) -------------
# Columns A and B carry explicit styles (A = bold/bordered index style,
# B = date/time number format). Clone that formatting from an existing
# data row (row 6) down onto the freshly appended rows so the new cells
# match the rest of the table.
$ws.Range("A6:B6").Copy()
$ws.Range("A21:B27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Write the refreshed dataset into rows 6-27 ----------------------
# Columns: A Index | B Timestamp | C Motorcycle(S) | D PKW(S) |
#          E Small Transporter(S) | F Transporter(S) | G LKW(S) |
#          H Motorcycle(R) | I PKW(R) | J Small Transporter(R) |
#          K Transporter(R) | L LKW(R) | M Bicycle
$rows = @(
    @(4,45250.34027777778,1,72,0,37,1,0,41,0,2,2,23),
    @(5,45250.34722222222,2,82,0,38,0,0,51,0,1,0,11),
    @(6,45250.35416666666,4,96,0,56,0,0,33,0,9,0,13),
    @(7,45250.36111111111,3,67,1,41,0,0,35,0,3,0,12),
    @(8,45392.33333333334,15,54,0,33,4,0,47,0,2,0,13),
    @(9,45392.34027777778,13,81,3,30,0,0,40,0,1,2,14),
    @(10,45392.34722222222,8,71,0,28,2,0,49,0,4,2,16),
    @(11,45392.35416666666,6,99,0,37,0,0,69,0,2,0,7),
    @(12,45392.36111111111,14,101,7,48,1,0,66,0,6,1,14),
    @(13,45392.36805555555,11,96,0,44,0,0,48,0,3,2,16),
    @(14,45392.66666666666,6,112,0,23,0,0,55,0,3,0,7),
    @(15,45392.67361111111,6,95,0,13,3,0,45,0,1,1,11),
    @(16,45392.68055555555,17,98,0,29,0,1,51,0,2,3,8),
    @(17,45392.6875,8,91,0,22,0,2,54,0,0,0,2),
    @(18,45392.69444444445,8,125,0,25,0,0,47,0,3,2,2),
    @(19,45392.70138888889,2,150,0,32,3,0,48,1,1,2,3),
    @(20,45392.83333333334,0,61,0,15,0,0,21,0,1,0,8),
    @(21,45392.84027777778,3,66,0,15,1,0,19,0,0,0,3),
    @(22,45392.84722222222,2,51,0,19,0,0,37,0,3,0,2),
    @(23,45392.85416666666,1,58,0,22,0,0,21,0,0,2,8),
    @(24,45392.86111111111,5,63,0,17,1,0,18,0,0,1,4),
    @(25,45392.86805555555,3,58,0,17,0,0,22,0,0,0,6)
)

$startRow = 6
$numRows = $rows.Count
$numCols = 13

$arr = New-Object 'object[,]' $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    $line = $rows[$r]
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r,$c] = $line[$c]
    }
}

$topLeft = $ws.Cells.Item($startRow, 1)
$bottomRight = $ws.Cells.Item($startRow + $numRows - 1, $numCols)
$target = $ws.Range($topLeft, $bottomRight)
$target.Value = $arr
